$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-18 04:47:03"
$wsDeDe.Range("H2").Value = "2016-08-18 04:47:03"
$wsZhCn.Range("H2").Value = "2016-08-18 04:46:56"
$wsZhCn.Range("K2").Value = "2016-08-18 04:47:28"
$wsDeDe.Range("K2").Value = "2016-08-18 04:47:35"
